$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D, shifting existing D:AM content (and the implicit
# empty trailing cells) one column to the right across the whole sheet.
# This reproduces: old D -> new E, old E -> new F, ... old AM -> new AN,
# leaving the new column D empty everywhere (rows 1-29).
$ws.Range("D1").EntireColumn.Insert()

# Row 23 header gets a brand-new first entry in the inserted column D.
$ws.Range("D23").Value = "rdf:type"

# Row 21: refreshed "modified" timestamp.
$ws.Range("B21").Value = "2023-08-17T09:08:29+00:00"

# New row 30: vocab:1006 / hasMaturityLevel / owl:ObjectProperty
$ws.Range("A30").Value = "vocab:1006"
$ws.Range("B30").Value = "hasMaturityLevel"
$ws.Range("D30").Value = "owl:ObjectProperty"
